# Daily scrape update - 2025-11-30 03:38:16 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh data rows 2-5 with the newly scraped opportunities ---
# (Opportunity IDs are kept as text via a leading apostrophe, matching
# how the source sheet stores them.)

# Row 2
$ws.Range("A2").Value = "'1329591"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1329591"
$ws.Range("C2").Value = "Research on Japan’s Energy Issues and Their Application to Another Nation"
$ws.Range("D2").Value = "日本、大阪府大阪市"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("H2").Value = "ITEM Corporation"

# Row 3
$ws.Range("A3").Value = "'1327958"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1327958"
$ws.Range("C3").Value = "Graphic designer"
$ws.Range("D3").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("F3").Value = "5 applicants"
$ws.Range("G3").Value = "3 - 6 Months"
$ws.Range("H3").Value = "The Paddock"

# Row 4
$ws.Range("A4").Value = "'1327495"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1327495"
$ws.Range("C4").Value = "Content Creator"
$ws.Range("D4").Value = "Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt"
$ws.Range("F4").Value = "9 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "Skyline Egypt Tours"

# Row 5
$ws.Range("A5").Value = "'1322114"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1322114"
$ws.Range("C5").Value = "Digital Marketing Specialist"
$ws.Range("D5").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("F5").Value = "30 applicants"
$ws.Range("G5").Value = "3 - 6 Months"
$ws.Range("H5").Value = "Fathalla & co"

# --- Drop the now-stale rows 6-9 (today's scrape only has 4 listings) ---
$ws.Rows("6:9").Delete()

# --- Widen columns C, D, G, H to fit the new content ---
# ColumnWidth is in "characters"; Excel persists the stored <col width="...">
# as characters + 5/6, so back that offset out to land exactly on the
# target stored widths (76, 56, 16, 22).
$ws.Columns("C").ColumnWidth = 76 - 5/6
$ws.Columns("D").ColumnWidth = 56 - 5/6
$ws.Columns("G").ColumnWidth = 16 - 5/6
$ws.Columns("H").ColumnWidth = 22 - 5/6
